$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# theta_se row (row 4): was all "(nan)", now pickled per-spec SEs
$ws.Range("B4").Value = "(0.0)"
$ws.Range("C4").Value = "(0.23)"
$ws.Range("D4").Value = "(0.08)"
$ws.Range("E4").Value = "(0.28)"
$ws.Range("F4").Value = "(0.22)"
$ws.Range("G4").Value = "(0.7)"

# lambda_se row (row 6): was all "(nan)", now pickled per-spec SEs
$ws.Range("B6").Value = "(0.0)"
$ws.Range("C6").Value = "(0.26)"
$ws.Range("D6").Value = "(0.17)"
$ws.Range("E6").Value = "(0.11)"
$ws.Range("F6").Value = "(0.02)"
$ws.Range("G6").Value = "(0.53)"
